# Actualización desde MV -datos-
# Updates quarterly financial-operations figures for rows 41-45 (revised
# figures for the existing quarters) and appends a new quarter row (46,
# "01-04-2021") at the bottom of the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Revisions to previously published quarters (rows 41-45)
# ---------------------------------------------------------------------

# Row 41 ("01-01-2020")
$ws.Cells.Item(41, 2).Value  = 805580     # B41
$ws.Cells.Item(41, 3).Value  = -1004074   # C41
$ws.Cells.Item(41, 5).Value  = 1181734    # E41
$ws.Cells.Item(41, 10).Value = 907042     # J41
$ws.Cells.Item(41, 16).Value = 958728     # P41
$ws.Cells.Item(41, 22).Value = -1468057   # V41
$ws.Cells.Item(41, 26).Value = 1469584    # Z41
$ws.Cells.Item(41, 28).Value = -153147    # AB41

# Row 42
$ws.Cells.Item(42, 2).Value  = 125611     # B42
$ws.Cells.Item(42, 3).Value  = 37036      # C42
$ws.Cells.Item(42, 5).Value  = 214561     # E42
$ws.Cells.Item(42, 10).Value = 173069     # J42
$ws.Cells.Item(42, 16).Value = 4653378    # P42
$ws.Cells.Item(42, 22).Value = 3118344    # V42
$ws.Cells.Item(42, 26).Value = 3225739    # Z42
$ws.Cells.Item(42, 28).Value = -4527767   # AB42

# Row 43
$ws.Cells.Item(43, 2).Value  = -4854321   # B43
$ws.Cells.Item(43, 3).Value  = 735754     # C43
$ws.Cells.Item(43, 5).Value  = 202676     # E43
$ws.Cells.Item(43, 10).Value = -182978    # J43
$ws.Cells.Item(43, 28).Value = -6168071   # AB43

# Row 44
$ws.Cells.Item(44, 2).Value  = -2265355   # B44
$ws.Cells.Item(44, 3).Value  = 429461     # C44
$ws.Cells.Item(44, 5).Value  = 391088     # E44
$ws.Cells.Item(44, 10).Value = -1133184   # J44
$ws.Cells.Item(44, 16).Value = 1161074    # P44
$ws.Cells.Item(44, 22).Value = 991064     # V44
$ws.Cells.Item(44, 26).Value = 1469611    # Z44
$ws.Cells.Item(44, 28).Value = -3426429   # AB44

# Row 45 ("01-01-2021")
$ws.Cells.Item(45, 2).Value  = -461109    # B45
$ws.Cells.Item(45, 3).Value  = -1153135   # C45
$ws.Cells.Item(45, 4).Value  = 183631     # D45
$ws.Cells.Item(45, 5).Value  = 1336766    # E45
$ws.Cells.Item(45, 10).Value = 1451569    # J45
$ws.Cells.Item(45, 16).Value = 633621     # P45
$ws.Cells.Item(45, 22).Value = -2385827   # V45
$ws.Cells.Item(45, 23).Value = 847461     # W45
$ws.Cells.Item(45, 26).Value = 3233288    # Z45
$ws.Cells.Item(45, 27).Value = -79477     # AA45
$ws.Cells.Item(45, 28).Value = -1094730   # AB45

# ---------------------------------------------------------------------
# New quarter row 46, "01-04-2021"
# ---------------------------------------------------------------------
# The label in column A must land in the workbook as plain shared-string
# text (matching every other "Serie" cell in column A), not as an
# auto-converted date serial. Typing a dd-mm-yyyy-look-alike string
# straight into .Value makes the COM layer parse it as a date (adding an
# unwanted number format/style). To avoid that, the label is built as a
# text formula in a scratch cell (forcing a text result), then brought
# into place with a values-only paste, which carries over the shared
# string but none of the scratch cell's formatting/style.
$scratchRow = 200
$scratchCol = 50
$scratch = $ws.Cells.Item($scratchRow, $scratchCol)
$scratch.Formula = '="01-04-2021"'
$scratch.Copy()
$ws.Cells.Item(46, 1).PasteSpecial(-4163)   # xlPasteValues
$scratch.Clear()

$ws.Cells.Item(46, 2).Value  = -739513    # B46
$ws.Cells.Item(46, 3).Value  = 598500     # C46
$ws.Cells.Item(46, 4).Value  = 971877     # D46
$ws.Cells.Item(46, 5).Value  = 373377     # E46
$ws.Cells.Item(46, 6).Value  = -2079742   # F46
$ws.Cells.Item(46, 7).Value  = 2982398    # G46
$ws.Cells.Item(46, 8).Value  = 5062141    # H46
$ws.Cells.Item(46, 9).Value  = 2987       # I46
$ws.Cells.Item(46, 10).Value = 738742     # J46
$ws.Cells.Item(46, 11).Value = 0          # K46
$ws.Cells.Item(46, 12).Value = 0          # L46
$ws.Cells.Item(46, 13).Value = 0          # M46
$ws.Cells.Item(46, 14).Value = 0          # N46
$ws.Cells.Item(46, 15).Value = 0          # O46
$ws.Cells.Item(46, 16).Value = 3522788    # P46
$ws.Cells.Item(46, 17).Value = 2435236    # Q46
$ws.Cells.Item(46, 18).Value = 2454077    # R46
$ws.Cells.Item(46, 21).Value = 18841      # U46
$ws.Cells.Item(46, 22).Value = 1164892    # V46
$ws.Cells.Item(46, 23).Value = 2518003    # W46
$ws.Cells.Item(46, 26).Value = 1353111    # Z46
$ws.Cells.Item(46, 27).Value = -77340     # AA46
$ws.Cells.Item(46, 28).Value = -4262301   # AB46
